$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("工作表1") ---
# Insert a new blank row at position 4. This pushes the old row 4 (A4=4 ...)
# down to row 5, leaving row 4 empty. Rows 2 and 3 (A2=2, A3=3 ...) stay put.
$ws1.Rows("4").Insert()

# The row that is now at row 5 used to hold A=4; move it back up to row 3,
# and put the row that used to hold A=3 down at row 5 (re-enter the values /
# formulas so everything recalculates with its new row number).
$ws1.Range("A3").Value = 4
$ws1.Range("B3").Formula = "=B`$1&`$A3"
$ws1.Range("C3").Formula = "=C`$1&`$A3"
$ws1.Range("D3").Formula = "=D`$1&`$A3"

$ws1.Range("A5").Value = 3
$ws1.Range("B5").Formula = "=B`$1&`$A5"
$ws1.Range("C5").Formula = "=C`$1&`$A5"
$ws1.Range("D5").Formula = "=D`$1&`$A5"

# New cells in columns E/F.
$ws1.Range("E2").Value = 2
$ws1.Range("E3").Value = 4
$ws1.Range("F3").Value = "iii"
$ws1.Range("E4").Value = "aaa"
$ws1.Range("F4").Value = "iii"
$ws1.Range("E5").Value = 3

# Page setup (paper size / orientation) now specified explicitly.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Sheet2 ("工作表1 (2)") ---
# Swap the "c" / "d" header labels in C1 / D1 (downstream formulas
# recalculate automatically since they reference the header row).
$ws2.Range("C1").Value = "d"
$ws2.Range("D1").Value = "c"

# New column E with header + values.
$ws2.Range("E1").Value = "b"
$ws2.Range("E2").Value = "zzz"
$ws2.Range("E3").Value = "xxx"

# Page setup (paper size / orientation) now specified explicitly.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- Selection / active sheet ---
# Sheet1's selection moves to C8; sheet2 becomes the active tab with E4
# selected (matches the tabSelected / activeTab change in the diff).
$null = $ws1.Range("C8").Select()
$null = $ws2.Select()
$null = $ws2.Range("E4").Select()
